$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$cell = $ws.Range("C1")
$cell.Value = 44307
$cell.NumberFormat = "mm-dd-yy"
